# Update generated-output numbers (view/like counts, min prices, image URL)
# for the 北京-漫展信息 workbook, matching the "output generated at 456a3b4"
# refresh described in the commit message / diff.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsShow    = $wb.Worksheets.Item("演出")
$wsLocal   = $wb.Worksheets.Item("本地生活")
$wsAll     = $wb.Worksheets.Item("全部类型")

# ---- 展览 (Exhibition) sheet ----
$wsExhibit.Range("F5").Value  = 360
$wsExhibit.Range("F6").Value  = 526
$wsExhibit.Range("F11").Value = 343
$wsExhibit.Range("F12").Value = 602
$wsExhibit.Range("F14").Value = 1495
$wsExhibit.Range("F15").Value = 1495
$wsExhibit.Range("F17").Value = 23
$wsExhibit.Range("F18").Value = 1343
$wsExhibit.Range("F20").Value = 273
$wsExhibit.Range("F24").Value = 6491
$wsExhibit.Range("F25").Value = 4823
$wsExhibit.Range("F26").Value = 134
$wsExhibit.Range("F28").Value = 204
$wsExhibit.Range("F29").Value = 138
$wsExhibit.Range("F32").Value = 1261
$wsExhibit.Range("F33").Value = 186
$wsExhibit.Range("F34").Value = 239
$wsExhibit.Range("F35").Value = 586
$wsExhibit.Range("F36").Value = 17
$wsExhibit.Range("F37").Value = 1333
$wsExhibit.Range("F38").Value = 231
$wsExhibit.Range("F40").Value = 141
$wsExhibit.Range("F43").Value = 90

# ---- 演出 (Performance) sheet ----
$wsShow.Range("F9").Value  = 2
$wsShow.Range("I9").Value  = "//i2.hdslb.com/bfs/openplatform/202408/B0Cuvd5v1724740500595.jpeg"
$wsShow.Range("G11").Value = 180
$wsShow.Range("F14").Value = 50

# ---- 本地生活 (Local life) sheet ----
$wsLocal.Range("F2").Value = 165
$wsLocal.Range("F3").Value = 2435
$wsLocal.Range("F4").Value = 179
$wsLocal.Range("F5").Value = 43

# ---- 全部类型 (All types) sheet ----
$wsAll.Range("F3").Value  = 165
$wsAll.Range("F7").Value  = 179
$wsAll.Range("F8").Value  = 43
$wsAll.Range("F9").Value  = 360
$wsAll.Range("F10").Value = 526
$wsAll.Range("F16").Value = 343
$wsAll.Range("F17").Value = 602
$wsAll.Range("F19").Value = 1495
$wsAll.Range("F20").Value = 1495
$wsAll.Range("F22").Value = 23
$wsAll.Range("F23").Value = 1343
$wsAll.Range("F25").Value = 273
$wsAll.Range("F30").Value = 6491
$wsAll.Range("F31").Value = 4823
$wsAll.Range("F32").Value = 134
$wsAll.Range("F33").Value = 204
$wsAll.Range("F35").Value = 1261
$wsAll.Range("F36").Value = 186
$wsAll.Range("F37").Value = 239
$wsAll.Range("G38").Value = 180
$wsAll.Range("F39").Value = 586
$wsAll.Range("F41").Value = 17
$wsAll.Range("F42").Value = 50
$wsAll.Range("F43").Value = 1333
$wsAll.Range("F44").Value = 231
$wsAll.Range("F45").Value = 141
$wsAll.Range("F48").Value = 90
